# Auto update stock data
# - Rename header F1 "Ticker" -> "Company"
# - Expand ticker symbols in column F to full company names (6-row blocks per company)
# - Refresh the first EBITDA value (column B) of each block with the latest figure
#
# Note: the EBITDA values in column B are stored as TEXT (not numbers) in this
# workbook. A plain `$ws.Range(...).Value = "4.59"` assignment would make Excel's
# smart cell-entry logic auto-convert a numeric-looking string into a real
# number, which would change the cell's stored type. To preserve the original
# text type (matching the source data), we stage the literal text in a scratch
# cell that has been explicitly formatted as Text, copy it, and use
# PasteSpecial to transplant the value (and its "keep as text" nature) onto the
# destination cell without disturbing the destination's own formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = "Z1"

function Set-TextValue([string]$cellRef, [string]$text) {
    $ws.Range($scratch).NumberFormat = "@"
    $ws.Range($scratch).Value = $text
    $ws.Range($scratch).Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)   # xlPasteValues
    $ws.Range($scratch).Clear()
}

# Column F header rename
$ws.Range("F1").Value = "Company"

# Each entry: ticker symbol (old), full company name (new), first row of its
# 6-row block, and the updated EBITDA value for that first row (column B).
# $null means no EBITDA update for that block.
$groups = @(
    @{ Ticker = "AA";      Company = "Alcoa";       StartRow = 2;  NewEbitda = "4.59" },
    @{ Ticker = "RIO";     Company = "Rio Tinto";    StartRow = 8;  NewEbitda = "7.51" },
    @{ Ticker = "NHY";     Company = "Norsk Hydro";  StartRow = 14; NewEbitda = "2.81" },
    @{ Ticker = "RS";      Company = "Reliance";     StartRow = 20; NewEbitda = $null },
    @{ Ticker = "KALU";    Company = "Kaiser";       StartRow = 26; NewEbitda = "9.82" },
    @{ Ticker = "RYI";     Company = "Ryerson";      StartRow = 32; NewEbitda = $null },
    @{ Ticker = "BVB:ALR"; Company = "Alro Steel";   StartRow = 38; NewEbitda = $null },
    @{ Ticker = "ULTR";    Company = "Ultra";        StartRow = 44; NewEbitda = "11.19" },
    @{ Ticker = "BHE";     Company = "Benchmark";    StartRow = 50; NewEbitda = $null },
    @{ Ticker = "CLS";     Company = "Celestica";    StartRow = 56; NewEbitda = $null },
    @{ Ticker = "JABIL";   Company = "Jabil";        StartRow = 62; NewEbitda = "11.55" },
    @{ Ticker = "FLEX";    Company = "Flex";         StartRow = 68; NewEbitda = "12.94" }
)

foreach ($g in $groups) {
    for ($i = 0; $i -lt 6; $i++) {
        $row = $g.StartRow + $i
        $ws.Range("F" + $row).Value = $g.Company
    }
    if ($g.NewEbitda -ne $null) {
        Set-TextValue ("B" + $g.StartRow) $g.NewEbitda
    }
}

# Final standalone EBITDA refresh (row 74, first row of the MKS block) -
# the ticker text for that block is unchanged in this update.
Set-TextValue "B74" "15.83"
